$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Google Cloud API를 활용하여 주소의 위도, 경도 정보 가져오기"
$ws.Range("E4").Value = "https://teddylee777.github.io/data_science/geocoding"

$ws.Range("D12").Value = "“머신러닝 파워드 애플리케이션” 사이킷런 1.0 업데이트 완료!"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/11/13/%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d-%ed%8c%8c%ec%9b%8c%eb%93%9c-%ec%95%a0%ed%94%8c%eb%a6%ac%ec%bc%80%ec%9d%b4%ec%85%98-%ec%82%ac%ec%9d%b4%ed%82%b7%eb%9f%b0-1-0-%ec%97%85%eb%8d%b0%ec%9d%b4/"

$ws.Range("D20").Value = "[머신러닝 기초] 비지도학습(Unsupervised-learning) - 군집화(Clustering)"
$ws.Range("E20").Value = "https://ai-creator.tistory.com/591"

$ws.Range("D44").Value = "메타버스를 위한 XR 비디오 표준 (MPEG-I/MIV)"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/107"
